$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (formatted YYYY-MM-DD).
# All rows from 2 to 110 currently hold 45181 and must be bumped to 45182.
$ws.Range("C2:C110").Value = 45182
